# Data update using git
# Update the "Inscritos" (column E) counts on the Inscricoes sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

$ws.Range("E2").Value = 30
$ws.Range("E12").Value = 33
$ws.Range("E14").Value = 40
$ws.Range("E15").Value = 107
$ws.Range("E16").Value = 327
